$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New skills: "Poison Status Effect" + fix for the OnTurnEnd trigger context bug ---
# Row 21 (ID 19): EfficientCasting / PassiveSkill
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "EfficientCasting"
$ws.Range("C21").Value = "PassiveSkill"
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0

# Row 22 (ID 20): PoisonMaw / DamageSkill(has effect)
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "PoisonMaw"
$ws.Range("C22").Value = "DamageSkill(has effect)"
$ws.Range("D22").Value = 10
$ws.Range("E22").Value = 3

# --- Normalize formatting: row 20 previously had mismatched styling
# (left over "Neutro"/"Ruim" look); make rows 20-22 match the uniform
# "Bom" (Good) look used by the rest of the data rows, by copying the
# format of an already-correct row onto them. ---
$src = $ws.Range("A19:E19")
$dst = $ws.Range("A20:E22")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the view to match where the user ended up after the edit ---
$ws.Range("A22:XFD22").Select()

"Done"
